# New crime data collected — weekly CompStat refresh (13th Precinct)
# Updates the report header (volume/issue number + week-covering dates)
# and the weekly/28-day/YTD/2-year crime figures for rows 15-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 32   Number  17" -> "...18" and the week-covering
# dates "4/21/2025"/"4/27/2025" -> "4/28/2025"/"5/4/2025".
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value = "Volume 32   Number  18"
$ws.Cells.Item(9, 3).Value = "Report Covering the Week  4/28/2025  Through  5/4/2025"

# ---------------------------------------------------------------------
# Helper donor cells used only to copy a *number format* (never touched
# by this edit) onto cells whose underlying type flips between a literal
# number and the sheet's placeholder text ("0" / "***.*"). Row 14 (the
# Murder row) is untouched by this week's data, so its cells are stable
# donors for style 13 (placeholder text), style 14 (integer count) and
# style 15 (percent-change figure).
# ---------------------------------------------------------------------
$styleTextDonor = $ws.Cells.Item(14, 3)   # C14 - style 13 (text / dash placeholder)
$styleCountDonor = $ws.Cells.Item(14, 10) # J14 - style 14 (integer count)
$stylePctDonor = $ws.Cells.Item(14, 11)   # K14 - style 15 (percent change)

function Set-NumberCell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-AsCount($row, $col, $value) {
    # Convert a placeholder-text cell into a literal number, reusing the
    # integer-count number format (style 14).
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
    $styleCountDonor.Copy()
    $c.PasteSpecial(-4122) # xlPasteFormats
}

function Set-AsPercent($row, $col, $value) {
    # Convert a placeholder-text cell into a literal number, reusing the
    # percent-change number format (style 15).
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
    $stylePctDonor.Copy()
    $c.PasteSpecial(-4122) # xlPasteFormats
}

function Set-AsPlaceholder($row, $col, $text) {
    # Convert a literal-number cell into the sheet's placeholder text
    # ("0" meaning "none reported" / "***.*" meaning "not computable"),
    # reusing the placeholder text's number format (style 13). The
    # leading apostrophe forces the numeric-looking "0" to be stored as
    # text instead of being re-parsed as a number.
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $styleTextDonor.Copy()
    $c.PasteSpecial(-4122) # xlPasteFormats
}

# ---------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------
Set-NumberCell 15 9 8
Set-NumberCell 15 11 300
Set-NumberCell 15 12 100
Set-NumberCell 15 13 60
Set-NumberCell 15 14 100

# ---------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------
Set-NumberCell 16 3 8
Set-NumberCell 16 5 166.666666666667
Set-NumberCell 16 6 18
Set-NumberCell 16 7 11
Set-NumberCell 16 8 63.636363636363
Set-NumberCell 16 9 57
Set-NumberCell 16 10 47
Set-NumberCell 16 11 21.276595744680
Set-NumberCell 16 12 3.636363636363
Set-NumberCell 16 13 -3.389830508474
Set-NumberCell 16 14 -82.515337423312

# ---------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------
Set-NumberCell 17 4 5
Set-NumberCell 17 5 0
Set-NumberCell 17 6 20
Set-NumberCell 17 7 16
Set-NumberCell 17 8 25
Set-NumberCell 17 9 82
Set-NumberCell 17 10 69
Set-NumberCell 17 11 18.840579710144
Set-NumberCell 17 12 15.492957746478
Set-NumberCell 17 13 127.777777777778
Set-NumberCell 17 14 -4.651162790697

# ---------------------------------------------------------------------
# Row 18 (C18 flips from a literal "5" to the "0" placeholder)
# ---------------------------------------------------------------------
Set-AsPlaceholder 18 3 "0"
Set-NumberCell 18 4 3
Set-NumberCell 18 5 -100
Set-NumberCell 18 6 15
Set-NumberCell 18 8 25
Set-NumberCell 18 9 128
Set-NumberCell 18 10 67
Set-NumberCell 18 11 91.044776119403
Set-NumberCell 18 12 47.126436781609
Set-NumberCell 18 13 26.732673267326
Set-NumberCell 18 14 -76.642335766423

# ---------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------
Set-NumberCell 19 3 30
Set-NumberCell 19 5 50
Set-NumberCell 19 6 90
Set-NumberCell 19 7 75
Set-NumberCell 19 8 20
Set-NumberCell 19 9 355
Set-NumberCell 19 10 332
Set-NumberCell 19 11 6.927710843373
Set-NumberCell 19 12 6.606606606606
Set-NumberCell 19 13 -26.041666666666
Set-NumberCell 19 14 -59.242250287026

# ---------------------------------------------------------------------
# Row 20 (D20/E20 flip from the "0"/"***.*" placeholders to literal
# numbers)
# ---------------------------------------------------------------------
Set-AsCount 20 4 1
Set-AsPercent 20 5 0
Set-NumberCell 20 6 5
Set-NumberCell 20 7 3
Set-NumberCell 20 8 66.666666666666
Set-NumberCell 20 9 12
Set-NumberCell 20 10 11
Set-NumberCell 20 11 9.090909090909
Set-NumberCell 20 12 -47.826086956521
Set-NumberCell 20 13 20
Set-NumberCell 20 14 -97.202797202797

# ---------------------------------------------------------------------
# Row 21 (bold "G.L.A." row, styles 17/18 unaffected by this edit)
# ---------------------------------------------------------------------
Set-NumberCell 21 3 45
Set-NumberCell 21 4 32
Set-NumberCell 21 5 40.625
Set-NumberCell 21 6 150
Set-NumberCell 21 7 117
Set-NumberCell 21 8 28.205128205128
Set-NumberCell 21 9 642
Set-NumberCell 21 10 529
Set-NumberCell 21 11 21.361058601134
Set-NumberCell 21 12 12.041884816753
Set-NumberCell 21 13 -7.091172214182
Set-NumberCell 21 14 -71.680635200705

# ---------------------------------------------------------------------
# Row 22 (F22 flips from a literal "2" to the "0" placeholder)
# ---------------------------------------------------------------------
Set-NumberCell 22 4 1
Set-AsPlaceholder 22 6 "0"
Set-NumberCell 22 8 -100
Set-NumberCell 22 10 23
Set-NumberCell 22 11 13.043478260869
Set-NumberCell 22 13 -10.344827586206

# ---------------------------------------------------------------------
# Row 23 (C23/F23 flip placeholder->number, D23/E23 flip number->placeholder)
# ---------------------------------------------------------------------
Set-AsCount 23 3 1
Set-AsPlaceholder 23 4 "0"
Set-AsPlaceholder 23 5 "***.*"
Set-AsCount 23 6 1
Set-NumberCell 23 8 0
Set-NumberCell 23 9 9
Set-NumberCell 23 11 200
Set-NumberCell 23 12 125
Set-NumberCell 23 13 50

# ---------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------
Set-NumberCell 24 3 42
Set-NumberCell 24 4 57
Set-NumberCell 24 5 -26.315789473684
Set-NumberCell 24 6 170
Set-NumberCell 24 7 238
Set-NumberCell 24 8 -28.571428571428
Set-NumberCell 24 9 877
Set-NumberCell 24 10 1060
Set-NumberCell 24 11 -17.264150943396
Set-NumberCell 24 12 28.781204111600
Set-NumberCell 24 13 49.914529914529

# ---------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------
Set-NumberCell 25 3 41
Set-NumberCell 25 4 49
Set-NumberCell 25 5 -16.326530612244
Set-NumberCell 25 6 149
Set-NumberCell 25 7 192
Set-NumberCell 25 8 -22.395833333333
Set-NumberCell 25 9 711
Set-NumberCell 25 10 898
Set-NumberCell 25 11 -20.824053452115
Set-NumberCell 25 12 38.596491228070

# ---------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------
Set-NumberCell 26 3 18
Set-NumberCell 26 4 11
Set-NumberCell 26 5 63.636363636363
Set-NumberCell 26 6 43
Set-NumberCell 26 7 42
Set-NumberCell 26 8 2.380952380952
Set-NumberCell 26 9 198
Set-NumberCell 26 10 186
Set-NumberCell 26 11 6.451612903225
Set-NumberCell 26 12 24.528301886792
Set-NumberCell 26 13 37.5

# ---------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------
Set-NumberCell 27 3 1
Set-NumberCell 27 9 12
Set-NumberCell 27 11 200
Set-NumberCell 27 12 33.333333333333

# ---------------------------------------------------------------------
# Row 28
# ---------------------------------------------------------------------
Set-NumberCell 28 3 6
Set-NumberCell 28 4 3
Set-NumberCell 28 5 100
Set-NumberCell 28 6 14
Set-NumberCell 28 7 13
Set-NumberCell 28 8 7.692307692307
Set-NumberCell 28 9 40
Set-NumberCell 28 10 42
Set-NumberCell 28 11 -4.761904761904
Set-NumberCell 28 12 -4.761904761904

# ---------------------------------------------------------------------
# Row 31
# ---------------------------------------------------------------------
Set-NumberCell 31 6 1
Set-NumberCell 31 12 150
